$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Remove the "Gott aber sei Dank ... (1. Korinther 15,57)" bible-verse row
#    (old row 8: category "Freude", bible ref "1 Kor 15,27"). Deleting the row
#    shifts every row below it up by one, which is what the target file shows
#    (dimension B2:D18 -> B2:D17, 14 data rows -> 13 data rows).
# ---------------------------------------------------------------------------
$ws.Rows.Item(8).Delete()

# ---------------------------------------------------------------------------
# 2. Re-stamp the sort range/condition so the persisted <sortState> matches
#    the new (one row shorter) table. The data is already sorted by the old
#    single-category column B, so re-applying the sort on the current range
#    is a no-op for row order but refreshes the stored range metadata to
#    B5:C17 / B5:B17.
# ---------------------------------------------------------------------------
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B5:B17"))
$ws.Sort.SetRange($ws.Range("B5:C17"))
$ws.Sort.Apply()

# ---------------------------------------------------------------------------
# 3. Update column B ("Kategorie/Motiv") on every data row to the new,
#    multi-category tags (the feature the commit message calls out: more
#    than one category/keyword per verse, including the new "Hand"/"Fluss"/
#    "Regenbogen" tags introduced alongside "Ansprache").
# ---------------------------------------------------------------------------
$newCategories = @(
    "Bild, Regenbogen",
    "Weg, Regenbogen, Bild",
    "Bild, Regenbogen",
    "Weg, Fluss, Psalm23, Hand",
    "Psalm23, Hand",
    "Psalm23, Weg, Fluss",
    "Psalm23, Hand, Säulen",
    "Psalm23, Weg",
    "Säulen, Hand",
    "Säulen",
    "Weg, Fluss, Regenbogen",
    "Weg, Fluss",
    "Weg, Fluss, Regenbogen"
)

for ($i = 0; $i -lt $newCategories.Length; $i++) {
    $row = 5 + $i
    $ws.Cells.Item($row, 2).Value = $newCategories[$i]
}

# ---------------------------------------------------------------------------
# 4. Re-point the two data-validation rules at their new target ranges. The
#    "custom" formula1 rule now covers B18:B25 and the "list" rule now only
#    covers B27 (both previously-scattered ranges collapse/shift once the
#    table above them is one row shorter).
# ---------------------------------------------------------------------------
$ws.Range("B23:B27").Validation.Delete()
$ws.Range("B5:B14").Validation.Delete()
$ws.Range("B16:B22").Validation.Delete()
$ws.Range("B29").Validation.Delete()

$ws.Range("B18:B25").Validation.Add(7, 1, 1, '"Weg;Freude;Psalm23;Säulen;Bild"')
$ws.Range("B27").Validation.Add(3, 1, 1, '"Weg,Freude,Psalm23,Säulen,Bild"')

# ---------------------------------------------------------------------------
# 5. Update the sheet view: the saved selection now sits at C19 (just below
#    the shortened table) and the view is scrolled so row 2 is the first
#    visible row.
# ---------------------------------------------------------------------------
$ws.Range("B2").Select()
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("C19").Select()

# ---------------------------------------------------------------------------
# 6. Minor column-width tweak on column B (Kategorie/Motiv) to fit the new,
#    longer multi-category labels.
# ---------------------------------------------------------------------------
$ws.Columns.Item(2).ColumnWidth = 22.1
